$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Files-tab Cypher query (cell B4) was corrected: the `File Type` and
# `Breed` lines were removed from the RETURN clause.
$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Boxer'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesQuery

# Row 4 shrinks now that the query text is two lines shorter.
$ws.Rows.Item(4).RowHeight = 217.5

# The author's selection ended up on B4 (scrolled so row 4 is at the top).
$ws.Range("B4").Select()
